$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 297
$ws1.Range("F6").Value = 273
$ws1.Range("F7").Value = 123
$ws1.Range("F12").Value = 122
$ws1.Range("F13").Value = 2557
$ws1.Range("F14").Value = 92
$ws1.Range("F17").Value = 21
$ws1.Range("F19").Value = 544
$ws1.Range("F20").Value = 611
$ws1.Range("F21").Value = 183
$ws1.Range("F22").Value = 95
$ws1.Range("F24").Value = 14
$ws1.Range("F26").Value = 2168
$ws1.Range("F27").Value = 4253
$ws1.Range("F29").Value = 67
$ws1.Range("F30").Value = 468
$ws1.Range("F31").Value = 1237
$ws1.Range("F32").Value = 248
$ws1.Range("F33").Value = 2143
$ws1.Range("F34").Value = 569
$ws1.Range("F35").Value = 475
$ws1.Range("F37").Value = 38
$ws1.Range("F38").Value = 136
$ws1.Range("F39").Value = 300
$ws1.Range("F41").Value = 737
$ws1.Range("F43").Value = 445
$ws1.Range("F44").Value = 20
$ws1.Range("F45").Value = 437

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 297
$ws4.Range("F6").Value = 273
$ws4.Range("F7").Value = 123
$ws4.Range("F12").Value = 122
$ws4.Range("F13").Value = 2557
$ws4.Range("F14").Value = 92
$ws4.Range("F18").Value = 21
$ws4.Range("F20").Value = 544
$ws4.Range("F21").Value = 611
$ws4.Range("F22").Value = 183
$ws4.Range("F23").Value = 95
$ws4.Range("F25").Value = 14
$ws4.Range("F27").Value = 2168
$ws4.Range("F28").Value = 4253
$ws4.Range("F30").Value = 67
$ws4.Range("F31").Value = 468
$ws4.Range("F32").Value = 1237
$ws4.Range("F33").Value = 248
$ws4.Range("F34").Value = 2143
$ws4.Range("F35").Value = 569
$ws4.Range("F36").Value = 475
$ws4.Range("F38").Value = 38
$ws4.Range("F39").Value = 136
$ws4.Range("F40").Value = 300
$ws4.Range("F42").Value = 737
$ws4.Range("F44").Value = 445
$ws4.Range("F45").Value = 20
$ws4.Range("F46").Value = 437
